# Append a new row (row 39) of sensor data to each of the 4 worksheets.
$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = 1
        A = "2025-03-05 22:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = 2
        A = "2025-03-05 22:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = 3
        A = "2025-03-05 22:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = 4
        A = "2025-03-05 22:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($r in $rows) {
    $ws = $wb.Worksheets.Item($r.Sheet)
    $rowIdx = 39

    $ws.Cells.Item($rowIdx, 1).Value = $r.A
    $ws.Cells.Item($rowIdx, 2).Value = $r.B
    $ws.Cells.Item($rowIdx, 3).Value = $r.C
    $ws.Cells.Item($rowIdx, 4).Value = $r.D
    $ws.Cells.Item($rowIdx, 5).Value = $r.E
    $ws.Cells.Item($rowIdx, 6).Value = $r.F

    # Column G holds a 24-digit identifier that exceeds double precision;
    # format as text first so Excel stores it verbatim instead of rounding
    # it into scientific notation.
    $gCell = $ws.Cells.Item($rowIdx, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $r.G

    $ws.Cells.Item($rowIdx, 8).Value = $r.H
    $ws.Cells.Item($rowIdx, 9).Value = $r.I
}
